$wb = $excel.ActiveWorkbook

$sheets = $wb.Worksheets

# Build the header row + style once on the first worksheet.
$ws1 = $sheets.Item(1)
$ws1.Range("A1").Value = "Input Sheet"
$ws1.Range("B1").Value = "Value"

$styleSource = $ws1.Range("A1")
$styleSource.Font.Bold = $true
$styleSource.Borders.LineStyle = 1
$styleSource.HorizontalAlignment = -4108
$styleSource.VerticalAlignment = -4160

$styleSource.Copy()
$ws1.Range("B1").PasteSpecial(-4122)

# Replicate the same header + style to every other worksheet by pasting
# the already-built formatting (avoids re-deriving new style records).
for ($i = 2; $i -le $sheets.Count; $i++) {
    $ws = $sheets.Item($i)
    $ws.Range("A1").Value = "Input Sheet"
    $ws.Range("B1").Value = "Value"

    $styleSource.Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("B1").PasteSpecial(-4122)
}
